# Applies the "12th day first commit" update to the TaskList sheet:
# adds 8 new task rows (106-113) continuing the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")

$date = Get-Date -Year 2017 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0

# Duplicate the formatting of the last two existing rows (104/105, which
# alternate the banded border style) down across the eight new rows so the
# new cells pick up the same styles already used throughout the sheet.
$ws.Range("A104:H105").Copy($ws.Range("A106:H107"))
$ws.Range("A104:H105").Copy($ws.Range("A108:H109"))
$ws.Range("A104:H105").Copy($ws.Range("A110:H111"))
$ws.Range("A104:H105").Copy($ws.Range("A112:H113"))
$excel.CutCopyMode = 0

$rows = @(
    @{ SNo = 105; Task = "Single Product Page";                 Time = "180 Minutes" },
    @{ SNo = 106; Task = "ProductCRUD Delete operation";         Time = "60Minutes" },
    @{ SNo = 107; Task = "Cart View Page";                       Time = "60Minutes" },
    @{ SNo = 108; Task = "Cart DTO, DAO, DAOIMPL";               Time = "60Minutes" },
    @{ SNo = 109; Task = "CartItem DTO, DAO, DAOIMPL";           Time = "60Minutes" },
    @{ SNo = 110; Task = "Payment view page";                    Time = "40Minutes" },
    @{ SNo = 111; Task = "Payment DTO";                          Time = "120Minutes" },
    @{ SNo = 112; Task = "CartFlow, CartModel, CartHandler";     Time = "120Minutes" }
)

$startRow = 106
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.SNo
    $ws.Cells.Item($r, 2).Value = $date
    $ws.Cells.Item($r, 3).Value = $data.Task
    $ws.Cells.Item($r, 6).Value = $data.Time
}

# Row 106 wraps to two lines ("180 Minutes" in the narrow Time Taken
# column), so its row grows to double height, matching what Excel would
# compute automatically once the text is entered.
$ws.Rows.Item(106).RowHeight = 28.8

# Update the view to match the saved state (scrolled down, selection on D110)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("D110").Select()
